$d = $word.ActiveDocument

$replacements = @(
    @("2024-06-21 Friday", "2024-06-22 Saturday"),
    @("562×6=3372", "870×5=4350"),
    @("139×6=834", "278×6=1668"),
    @("872×9=7848", "663×6=3978"),
    @("423×6=2538", "614×9=5526"),
    @("739×5=3695", "397×9=3573"),
    @("880×6=5280", "933×9=8397"),
    @("183×4=732", "292×4=1168"),
    @("294×9=2646", "780×6=4680"),
    @("755×6=4530", "838×8=6704"),
    @("837×2=1674", "644×6=3864"),
    @("663×4=2652", "230×5=1150"),
    @("634×8=5072", "431×9=3879"),
    @("611×3=1833", "415×7=2905"),
    @("828×4=3312", "654×3=1962"),
    @("958×2=1916", "359×5=1795"),
    @("296×4=1184", "521×6=3126"),
    @("165×8=1320", "963×7=6741"),
    @("719×5=3595", "317×7=2219"),
    @("551×9=4959", "504×2=1008"),
    @("664×7=4648", "991×2=1982"),
    @("347×8=2776", "714×9=6426"),
    @("649×7=4543", "819×8=6552"),
    @("569×5=2845", "581×6=3486"),
    @("240×3=720", "522×9=4698"),
    @("313×2=626", "609×6=3654")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
